$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column layout change: the "Mẫu 9" template used to have 14 columns
#    (A:N) with separate "Biển số" / "Màu biển" columns and two
#    "Ngày cấp CCCD..." columns. The new layout merges "Biển số" and
#    "Màu biển" into a single "Biển số/ màu biển" column and removes the
#    two "Ngày cấp CCCD..." columns entirely, leaving 11 columns (A:K).
#
#    Deleting old column C ("Màu biển") folds it into column B ("Biển số").
#    After that delete, the old "Ngày cấp CCCD của chủ xe" column has
#    shifted left to H, and the old "Ngày cấp CCCD/mã số thuế người mua"
#    column has shifted left to K - delete both of those too.
# ---------------------------------------------------------------------------
$ws.Columns("C").Delete()
$ws.Columns("H").Delete()
$ws.Columns("K").Delete()

# Rename the merged "Biển số" header (now column B) to reflect both
# the plate number and plate colour.
$ws.Range("B4").Value = "Biển số/ màu biển"

# ---------------------------------------------------------------------------
# 2. The second section of the sheet (list continuation starting at row 18)
#    used to repeat only the "STT" header on row 17. Now it repeats the
#    full header row, same as row 4.
# ---------------------------------------------------------------------------
$ws.Range("B4:K4").Copy()
$ws.Range("B17:K17").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row heights were nudged slightly (minor re-layout / re-save), and row
#    17 now needs an explicit height to fit its newly-populated, wrapped
#    header text.
# ---------------------------------------------------------------------------
$ws.Rows("1").RowHeight = 58.2
$ws.Rows("2").RowHeight = 13.8
$ws.Rows("4").RowHeight = 111.6
$ws.Rows("5:16").RowHeight = 28.05
$ws.Rows("17").RowHeight = 105.6
$ws.Rows("18:32").RowHeight = 28.05

# ---------------------------------------------------------------------------
# 4. Update the active selection to match where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("D29").Select()
